$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new translation rows (status "new", Khmer not yet translated -> "TBT"),
# each inserted at the correct alphabetical position. Rows are inserted top-to-bottom;
# because each insertion shifts subsequent rows down by one, the row indices below are
# the live (current) row numbers at the moment of each insertion.

# 1) "Clinical data not provided" -> before "Clinical Outcome" (row 33)
$ws.Rows.Item(33).Insert()
$ws.Cells.Item(33,1).Value2 = "Clinical data not provided"
$ws.Cells.Item(33,2).Value2 = "TBT"
$ws.Cells.Item(33,3).Value2 = "new"

# 2) "Lab data not provided" -> before "Lab data successfully processed!" (row 80 after shift above)
$ws.Rows.Item(80).Insert()
$ws.Cells.Item(80,1).Value2 = "Lab data not provided"
$ws.Cells.Item(80,2).Value2 = "TBT"
$ws.Cells.Item(80,3).Value2 = "new"

# 3) "No .acorn has been generated" -> before "No Blood Culture" (row 95 after shifts above)
$ws.Rows.Item(95).Insert()
$ws.Cells.Item(95,1).Value2 = "No .acorn has been generated"
$ws.Cells.Item(95,2).Value2 = "TBT"
$ws.Cells.Item(95,3).Value2 = "new"

# 4) "No .acorn has been saved" -> right after the row above, still before "No Blood Culture" (row 96)
$ws.Rows.Item(96).Insert()
$ws.Cells.Item(96,1).Value2 = "No .acorn has been saved"
$ws.Cells.Item(96,2).Value2 = "TBT"
$ws.Cells.Item(96,3).Value2 = "new"
